# Simplify the "conta" labels in column B: the long, versioned Portuguese
# descriptions are replaced with their short canonical names. The "exercicio"
# dates in column A and the "valor" numbers in column C are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10  -> "RECEITA CORRENTE LÍQUIDA - RCL" / "... RCL (IV)"
$ws.Range("B2:B10").Value = "RECEITA CORRENTE LÍQUIDA"

# Rows 11-19 -> "DÍVIDA CONSOLIDADA LÍQUIDA (DCL) (III) = (I - II)"
$ws.Range("B11:B19").Value = "DÍVIDA CONSOLIDADA LÍQUIDA"

# Rows 20-28 -> "% da DCL sobre a RCL (III/RCL)" / "... AJUSTADA (III/VI)"
$ws.Range("B20:B28").Value = "% da DCL sobre a RCL"

# Column B now holds much longer text than before ("DÍVIDA CONSOLIDADA
# LÍQUIDA" is the longest), so widen it to fit, matching the workbook's
# existing best-fit styling on column A.
$ws.Columns("B").ColumnWidth = 44
